$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting the existing row 22 (and below) down to row 23.
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the new data record.
$ws.Cells.Item(22, 1).Value = 11
$ws.Cells.Item(22, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(22, 3).Value = "Bíobío"
$ws.Cells.Item(22, 4).Value = 44595
$ws.Cells.Item(22, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(22, 6).Value = 100112022
$ws.Cells.Item(22, 7).Value = "Arveja Verde"
$ws.Cells.Item(22, 8).Value = "Perfection"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 50
$ws.Cells.Item(22, 11).Value = 26000
$ws.Cells.Item(22, 12).Value = 28000
$ws.Cells.Item(22, 13).Value = 27200
$ws.Cells.Item(22, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(22, 15).Value = "Carahue"
$ws.Cells.Item(22, 16).Value = 1088
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"
